$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "'1"
$ws.Range("D2").Value = "'3"

$ws.Range("C3").Value = "'0"
$ws.Range("D3").Value = "'1"
$ws.Range("E3").Value = "'0"

$ws.Range("C4").Value = "'9"
$ws.Range("D4").Value = "'10"
$ws.Range("E4").Value = "'1"
